# chore: update pv pvc volumes
#
# The "variables" worksheet holds one VARIABLE/AKS/SUNAT triple per row
# (columns A/B/C), with column D a helper formula that concatenates the
# three into a single "VARIABLE:x|AKS:y|SUNAT:z" string. This change adds
# one new row for the "volumeSize" variable (AKS value 8Gi, SUNAT value
# 90Gi) right after the last existing row (50 - csiDriverNfs).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 51

# Write B/C before A so new shared-string entries land in the same order
# the source workbook uses (8Gi, 90Gi, volumeSize).
$ws.Cells.Item($newRow, 2).Value = "8Gi"
$ws.Cells.Item($newRow, 3).Value = "90Gi"
$ws.Cells.Item($newRow, 1).Value = "volumeSize"

# Columns B and C are stored as text (same as every other data row) -
# match that formatting for the new cells.
$ws.Cells.Item($newRow, 2).NumberFormat = "@"
$ws.Cells.Item($newRow, 3).NumberFormat = "@"

# Column D repeats the same concatenation formula used by every row above.
$ws.Cells.Item($newRow, 4).Formula = '=$A$1&":"&A51&"|"&$B$1&":"&B51&"|"&$C$1&":"&C51'

# Match the author's final selection: cell A51 (the first cell of the new
# row) becomes the active/selected cell.
$ws.Range("A51").Select()
